$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Updated "time_taken" timestamps for rows 2-105 (column F), reflecting a
# newer panel query run time.
$newTimes = @(
    "2021-10-05 14:19:30.170121",
    "2021-10-05 14:19:30.170129",
    "2021-10-05 14:19:30.170132",
    "2021-10-05 14:19:30.170135",
    "2021-10-05 14:19:30.170137",
    "2021-10-05 14:19:30.170140",
    "2021-10-05 14:19:30.170142",
    "2021-10-05 14:19:30.170145",
    "2021-10-05 14:19:30.170148",
    "2021-10-05 14:19:30.170151",
    "2021-10-05 14:19:30.170153",
    "2021-10-05 14:19:30.170156",
    "2021-10-05 14:19:30.170158",
    "2021-10-05 14:19:30.170161",
    "2021-10-05 14:19:30.170164",
    "2021-10-05 14:19:30.170166",
    "2021-10-05 14:19:30.170169",
    "2021-10-05 14:19:30.170172",
    "2021-10-05 14:19:30.170175",
    "2021-10-05 14:19:30.170177",
    "2021-10-05 14:19:30.170180",
    "2021-10-05 14:19:30.170182",
    "2021-10-05 14:19:30.170185",
    "2021-10-05 14:19:30.170187",
    "2021-10-05 14:19:30.170190",
    "2021-10-05 14:19:30.170193",
    "2021-10-05 14:19:30.170195",
    "2021-10-05 14:19:30.170198",
    "2021-10-05 14:19:30.170200",
    "2021-10-05 14:19:30.170203",
    "2021-10-05 14:19:30.170205",
    "2021-10-05 14:19:30.170208",
    "2021-10-05 14:19:30.170211",
    "2021-10-05 14:19:30.170213",
    "2021-10-05 14:19:30.170216",
    "2021-10-05 14:19:30.170218",
    "2021-10-05 14:19:30.170221",
    "2021-10-05 14:19:30.170223",
    "2021-10-05 14:19:30.170226",
    "2021-10-05 14:19:30.170228",
    "2021-10-05 14:19:30.170232",
    "2021-10-05 14:19:30.170234",
    "2021-10-05 14:19:30.170237",
    "2021-10-05 14:19:30.170240",
    "2021-10-05 14:19:30.170242",
    "2021-10-05 14:19:30.170245",
    "2021-10-05 14:19:30.170247",
    "2021-10-05 14:19:30.170250",
    "2021-10-05 14:19:30.170252",
    "2021-10-05 14:19:30.170255",
    "2021-10-05 14:19:30.170257",
    "2021-10-05 14:19:30.170260",
    "2021-10-05 14:19:30.170263",
    "2021-10-05 14:19:30.170265",
    "2021-10-05 14:19:30.170268",
    "2021-10-05 14:19:30.170271",
    "2021-10-05 14:19:30.170273",
    "2021-10-05 14:19:30.170276",
    "2021-10-05 14:19:30.170278",
    "2021-10-05 14:19:30.170281",
    "2021-10-05 14:19:30.170284",
    "2021-10-05 14:19:30.170286",
    "2021-10-05 14:19:30.170289",
    "2021-10-05 14:19:30.170292",
    "2021-10-05 14:19:30.170295",
    "2021-10-05 14:19:30.170298",
    "2021-10-05 14:19:30.170301",
    "2021-10-05 14:19:30.170304",
    "2021-10-05 14:19:30.170307",
    "2021-10-05 14:19:30.170310",
    "2021-10-05 14:19:30.170312",
    "2021-10-05 14:19:30.170315",
    "2021-10-05 14:19:30.170318",
    "2021-10-05 14:19:30.170321",
    "2021-10-05 14:19:30.170323",
    "2021-10-05 14:19:30.170326",
    "2021-10-05 14:19:30.170330",
    "2021-10-05 14:19:30.170333",
    "2021-10-05 14:19:30.170336",
    "2021-10-05 14:19:30.170339",
    "2021-10-05 14:19:30.170341",
    "2021-10-05 14:19:30.170344",
    "2021-10-05 14:19:30.170347",
    "2021-10-05 14:19:30.170349",
    "2021-10-05 14:19:30.170352",
    "2021-10-05 14:19:30.170354",
    "2021-10-05 14:19:30.170357",
    "2021-10-05 14:19:30.170360",
    "2021-10-05 14:19:30.170362",
    "2021-10-05 14:19:30.170365",
    "2021-10-05 14:19:30.170368",
    "2021-10-05 14:19:30.170370",
    "2021-10-05 14:19:30.170374",
    "2021-10-05 14:19:30.170377",
    "2021-10-05 14:19:30.170380",
    "2021-10-05 14:19:30.170383",
    "2021-10-05 14:19:30.170386",
    "2021-10-05 14:19:30.170388",
    "2021-10-05 14:19:30.170391",
    "2021-10-05 14:19:30.170393",
    "2021-10-05 14:19:30.170396",
    "2021-10-05 14:19:30.170399",
    "2021-10-05 14:19:30.170401",
    "2021-10-05 14:19:30.170404"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $newTimes[$i]
}

# Add a new "metadata" worksheet right after the "data" sheet.
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$metaSheet.Name = "metadata"

# Header row (bold, bordered, centered - same look as the "data" header row)
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"
$metaSheet.Range("B1:G1").Font.Bold = $true
$metaSheet.Range("B1:G1").Borders.LineStyle = 1
$metaSheet.Range("B1:G1").HorizontalAlignment = -4108
$metaSheet.Range("B1:G1").VerticalAlignment = -4160

# Data row
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Range("A2").Font.Bold = $true
$metaSheet.Range("A2").Borders.LineStyle = 1
$metaSheet.Range("A2").HorizontalAlignment = -4108
$metaSheet.Range("A2").VerticalAlignment = -4160

$metaSheet.Cells.Item(2, 2).Value = "Cerebral vascular malformations"
$metaSheet.Cells.Item(2, 3).Value = 147

$d2 = $metaSheet.Cells.Item(2, 4)
$d2.NumberFormat = "@"
$d2.Value = "2.58"

$metaSheet.Cells.Item(2, 5).Value = "2021-09-06T10:16:44.593839Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:19:30.166778"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/147/?format=json"

$metaSheet.Range("A1").Select() | Out-Null

# Keep "data" as the active sheet/tab, matching the original workbook view
# (only a new sheet entry was added; the active tab stays on "data").
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
